$d = $word.ActiveDocument

# Replace the document body text.
$d.Content.Find.Execute("Test Word Document", $true, $false, $false, $false, $false, $true, 1, $false,
                         "Call 555-1234 to buy ice cream, or the world will end!", 2)

# Apply character formatting (red, 48pt) to the whole paragraph, which also
# stamps the paragraph mark's run properties (w:pPr/w:rPr).
$para = $d.Paragraphs.Item(1)
$rng = $para.Range
$rng.Font.Color = 255
$rng.Font.Size = 48
$rng.Font.SizeBi = 48

# Touch the footnote/endnote machinery so Word mints word/footnotes.xml and
# word/endnotes.xml (with the standard separator/continuationSeparator
# boilerplate), then remove the scratch note itself so no reference mark is
# left behind in the body.
$endRng = $d.Content
$endRng.Collapse(0)
$fn = $d.Footnotes.Add($endRng, "", "x")
$fn.Delete()
